# Apply cryptos list update (price/volume refresh + two coin-pair swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.380.54"
$ws.Cells.Item(2, 5).Value = "  +3.17%  "

$ws.Cells.Item(3, 4).Value = "1.871.56"
$ws.Cells.Item(3, 5).Value = "  +1.47%  "

$ws.Cells.Item(4, 5).Value = "  -0.40%  "

$c = $ws.Cells.Item(5, 4)
$c.Value = "'338.97"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.55%  "

$ws.Cells.Item(6, 5).Value = "  -0.42%  "

$c = $ws.Cells.Item(7, 4)
$c.Value = "'0.4701"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +1.95%  "

$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.3960"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +3.24%  "

$c = $ws.Cells.Item(9, 4)
$c.Value = "'47.46"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +2.02%  "

$c = $ws.Cells.Item(10, 4)
$c.Value = "'0.08018"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +1.58%  "

$c = $ws.Cells.Item(11, 4)
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.98%  "

$c = $ws.Cells.Item(12, 4)
$c.Value = "'21.85"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +2.90%  "

$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.866.42"
$ws.Cells.Item(13, 5).Value = "  +0.96%  "

$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Cells.Item(14, 4)
$c.Value = "'5.989"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.33%  "

$c = $ws.Cells.Item(15, 4)
$c.Value = "'7.236"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +2.70%  "

$c = $ws.Cells.Item(16, 4)
$c.Value = "'91.21"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +3.65%  "

$ws.Cells.Item(17, 5).Value = "  -0.44%  "

$c = $ws.Cells.Item(18, 4)
$c.Value = "'0.00001042"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.10%  "

$c = $ws.Cells.Item(19, 4)
$c.Value = "'0.06615"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.31%  "

$c = $ws.Cells.Item(20, 4)
$c.Value = "'17.55"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +3.32%  "

$ws.Cells.Item(21, 5).Value = "  -0.45%  "

$ws.Cells.Item(22, 4).Value = "28.389.09"

$c = $ws.Cells.Item(23, 4)
$c.Value = "'5.454"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.84%  "

$ws.Cells.Item(24, 5).Value = "  +1.37%  "

$c = $ws.Cells.Item(25, 4)
$c.Value = "'2.268"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -1.18%  "

$ws.Cells.Item(26, 4).Value = "2.088.28"
$ws.Cells.Item(26, 5).Value = "  +0.90%  "

$c = $ws.Cells.Item(27, 4)
$c.Value = "'160.23"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +1.97%  "

$c = $ws.Cells.Item(28, 4)
$c.Value = "'19.77"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.81%  "

$c = $ws.Cells.Item(29, 4)
$c.Value = "'2.130"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.49%  "

$c = $ws.Cells.Item(30, 4)
$c.Value = "'5.513"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +2.89%  "

$c = $ws.Cells.Item(31, 4)
$c.Value = "'120.02"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.91%  "

$c = $ws.Cells.Item(32, 4)
$c.Value = "'0.9693"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.10%  "

$c = $ws.Cells.Item(33, 4)
$c.Value = "'0.09482"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +1.57%  "

$c = $ws.Cells.Item(34, 4)
$c.Value = "'3.572"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.03%  "

$c = $ws.Cells.Item(35, 4)
$c.Value = "'1.377"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +3.95%  "

$c = $ws.Cells.Item(36, 4)
$c.Value = "'5.349"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +1.92%  "

$c = $ws.Cells.Item(37, 4)
$c.Value = "'0.06092"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +2.40%  "

$c = $ws.Cells.Item(38, 4)
$c.Value = "'0.02254"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +2.24%  "

$c = $ws.Cells.Item(39, 4)
$c.Value = "'8.374"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +3.34%  "

$c = $ws.Cells.Item(40, 4)
$c.Value = "'1.185"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.09%  "

$c = $ws.Cells.Item(41, 4)
$c.Value = "'0.5945"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +1.67%  "

$ws.Cells.Item(42, 5).Value = "  -0.48%  "

$c = $ws.Cells.Item(43, 4)
$c.Value = "'0.1871"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.29%  "

$c = $ws.Cells.Item(44, 4)
$c.Value = "'10.33"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.91%  "

$c = $ws.Cells.Item(45, 4)
$c.Value = "'1.293"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +2.90%  "

$c = $ws.Cells.Item(46, 4)
$c.Value = "'0.5588"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.13%  "

$c = $ws.Cells.Item(47, 4)
$c.Value = "'12.16"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.13%  "

$c = $ws.Cells.Item(48, 4)
$c.Value = "'1.956"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +4.55%  "

$c = $ws.Cells.Item(49, 4)
$c.Value = "'0.06865"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.08%  "

$ws.Cells.Item(50, 2).Value = "Quant"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Cells.Item(50, 4)
$c.Value = "'111.37"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.79%  "

$ws.Cells.Item(51, 2).Value = "RenderToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Cells.Item(51, 4)
$c.Value = "'2.035"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +13.05%  "
